$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 288; this pushes the existing rows 288-299
# down to 289-300, matching every per-row value shift seen in the diff.
$ws.Rows.Item(288).Insert()

# Populate the newly inserted row 288 with the new weekly price record.
$ws.Cells.Item(288, 1).Value = 7
$ws.Cells.Item(288, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(288, 3).Value = "Ñuble"
$ws.Cells.Item(288, 4).Value = 45147
$ws.Cells.Item(288, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(288, 5).Value = 16
$ws.Cells.Item(288, 6).Value = 100112040
$ws.Cells.Item(288, 7).Value = "Cilantro"
$ws.Cells.Item(288, 8).Value = "Sin especificar"
$ws.Cells.Item(288, 9).Value = "Primera"
$ws.Cells.Item(288, 10).Value = 180
$ws.Cells.Item(288, 11).Value = 1500
$ws.Cells.Item(288, 12).Value = 1500
$ws.Cells.Item(288, 13).Value = 1500
$ws.Cells.Item(288, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(288, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(288, 16).Value = 1500
$ws.Cells.Item(288, 17).Value = 1
$ws.Cells.Item(288, 18).Value = "Hortaliza"
